# Sheet "10per change" (first sheet) holds the breakout screener data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 35-37: bsecode (col E) was stored as text, now numeric ---
$ws.Cells.Item(35, 5).Value = 20
$ws.Cells.Item(36, 5).Value = 531344
$ws.Cells.Item(37, 5).Value = 505537

# --- New row 38 ---
$ws.Cells.Item(38, 1).Value = "24/06/2024 07:44:48"
$ws.Cells.Item(38, 2).Value = 1
$ws.Cells.Item(38, 3).Value = "BSE"
$ws.Cells.Item(38, 4).Value = "BSE (Bombay stock exchange)"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = "20"
$ws.Cells.Item(38, 5).Style = "Normal"
$ws.Cells.Item(38, 6).Value = -2.34
$ws.Cells.Item(38, 7).Value = 2500.25
$ws.Cells.Item(38, 8).Value = 669979

# --- New row 39 ---
$ws.Cells.Item(39, 1).Value = "24/06/2024 07:44:48"
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = "CONCOR"
$ws.Cells.Item(39, 4).Value = "Container Corporation Of India Limited"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = "531344"
$ws.Cells.Item(39, 5).Style = "Normal"
$ws.Cells.Item(39, 6).Value = -3.13
$ws.Cells.Item(39, 7).Value = 1056.8
$ws.Cells.Item(39, 8).Value = 3344893
